$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("C2").Value = 3
    $ws.Range("F2").Value = 0
    $ws.Range("H2").Value = 72
    $ws.Range("D3").Value = 34
    $ws.Range("F3").Value = 9
    $ws.Range("G3").Value = 8
    $ws.Range("H3").Value = 45
    $ws.Range("D4").Value = 36
    $ws.Range("E4").Value = 19
    $ws.Range("F4").Value = 3
    $ws.Range("G4").Value = 6
    $ws.Range("H4").Value = 43
    $ws.Range("D5").Value = 36
    $ws.Range("E5").Value = 16
    $ws.Range("F5").Value = 8
    $ws.Range("G5").Value = 10
    $ws.Range("H5").Value = 37
    $ws.Range("D6").Value = 37
    $ws.Range("E6").Value = 18
    $ws.Range("F6").Value = 8
    $ws.Range("G6").Value = 6
    $ws.Range("H6").Value = 38
    $ws.Range("D7").Value = 35
    $ws.Range("E7").Value = 13
    $ws.Range("F7").Value = 7
    $ws.Range("G7").Value = 12
    $ws.Range("H7").Value = 40
    $ws.Range("D8").Value = 36
    $ws.Range("E8").Value = 10
    $ws.Range("F8").Value = 10
    $ws.Range("G8").Value = 10
    $ws.Range("H8").Value = 41
    $ws.Range("D9").Value = 40
    $ws.Range("E9").Value = 15
    $ws.Range("F9").Value = 11
    $ws.Range("G9").Value = 6
    $ws.Range("H9").Value = 35
    $ws.Range("D10").Value = 29
    $ws.Range("E10").Value = 10
    $ws.Range("F10").Value = 6
    $ws.Range("G10").Value = 7
    $ws.Range("H10").Value = 55
    $ws.Range("D11").Value = 38
    $ws.Range("E11").Value = 10
    $ws.Range("F11").Value = 9
    $ws.Range("G11").Value = 7
    $ws.Range("H11").Value = 43
    $ws.Range("D12").Value = 32
    $ws.Range("E12").Value = 15
    $ws.Range("F12").Value = 3
    $ws.Range("G12").Value = 8
    $ws.Range("H12").Value = 49
    $ws.Range("E13").Value = 19
    $ws.Range("F13").Value = 7
    $ws.Range("G13").Value = 7
    $ws.Range("H13").Value = 37
    $ws.Range("D14").Value = 38
    $ws.Range("E14").Value = 18
    $ws.Range("F14").Value = 12
    $ws.Range("G14").Value = 10
    $ws.Range("H14").Value = 29
    $ws.Range("D15").Value = 35
    $ws.Range("E15").Value = 13
    $ws.Range("F15").Value = 8
    $ws.Range("G15").Value = 12
    $ws.Range("H15").Value = 39
    $ws.Range("D16").Value = 30
    $ws.Range("F16").Value = 5
    $ws.Range("G16").Value = 8
    $ws.Range("H16").Value = 55
    $ws.Range("D17").Value = 33
    $ws.Range("E17").Value = 15
    $ws.Range("H17").Value = 43
    $ws.Range("D18").Value = 44
    $ws.Range("E18").Value = 11
    $ws.Range("F18").Value = 7
    $ws.Range("G18").Value = 11
    $ws.Range("H18").Value = 34
    $ws.Range("D19").Value = 33
    $ws.Range("E19").Value = 6
    $ws.Range("H19").Value = 57
    $ws.Range("D20").Value = 33
    $ws.Range("E20").Value = 11
    $ws.Range("G20").Value = 7
    $ws.Range("H20").Value = 53
    $ws.Range("D21").Value = 37
    $ws.Range("E21").Value = 22
    $ws.Range("F21").Value = 8
    $ws.Range("G21").Value = 11
    $ws.Range("D22").Value = 37
    $ws.Range("E22").Value = 9
    $ws.Range("F22").Value = 7
    $ws.Range("G22").Value = 4
    $ws.Range("H22").Value = 50
    $ws.Range("D23").Value = 41
    $ws.Range("E23").Value = 14
    $ws.Range("F23").Value = 7
    $ws.Range("G23").Value = 9
    $ws.Range("H23").Value = 36
    $ws.Range("D24").Value = 45
    $ws.Range("E24").Value = 18
    $ws.Range("F24").Value = 10
    $ws.Range("G24").Value = 8
    $ws.Range("H24").Value = 26
    $ws.Range("D25").Value = 38
    $ws.Range("F25").Value = 9
    $ws.Range("G25").Value = 9
    $ws.Range("H25").Value = 37
    $ws.Range("D26").Value = 30
    $ws.Range("G26").Value = 10
    $ws.Range("H26").Value = 62
    $ws.Range("D27").Value = 35
    $ws.Range("E27").Value = 12
    $ws.Range("G27").Value = 8
    $ws.Range("H27").Value = 45
    $ws.Range("D28").Value = 33
    $ws.Range("E28").Value = 8
    $ws.Range("F28").Value = 5
    $ws.Range("G28").Value = 12
    $ws.Range("H28").Value = 49
    $ws.Range("E30").Value = 9
    $ws.Range("F30").Value = 5
    $ws.Range("H30").Value = 45
    $ws.Range("D31").Value = 37
    $ws.Range("E31").Value = 20
    $ws.Range("F31").Value = 9
    $ws.Range("G31").Value = 14
    $ws.Range("H31").Value = 27
    $ws.Range("D32").Value = 37
    $ws.Range("E32").Value = 15
    $ws.Range("F32").Value = 7
    $ws.Range("G32").Value = 7
    $ws.Range("H32").Value = 41
